$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 6 (shifts rows 6.. down by one, formulas auto-adjust)
$ws.Rows("6:6").Insert()

# Copy formatting down from row 5 into the freshly inserted row 6
$ws.Range("B5:C5").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)

# Populate the newly inserted row with the "Gecorreleerd" / "Nee" pair
$ws.Range("B6").Value = "Gecorreleerd"
$ws.Range("C6").Value = "Nee"

# Keep the conditional-formatting ranges anchored to the shifted data block
$ws.Range("G16:H35").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G17:H36"))
$ws.Range("J16:K35").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J17:K36"))
$ws.Range("O16:O35").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("O17:O36"))

# Restore the view: select D5, no frozen/scrolled topLeftCell
$ws.Range("D5").Select()
